$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1199.8334
$ws.Range("I100").Value = 1227.091
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1227.091
$ws.Range("L100").Value = 900
$ws.Range("M100").Value = -686.0909999999999
$ws.Range("N100").Value = -1982
$ws.Range("H127").Value = 829.65515
$ws.Range("J127").Value = 993.3684
$ws.Range("L127").Value = 2980.1052
$ws.Range("N127").Value = -12900.1052
$ws.Range("H129").Value = 947.0526
$ws.Range("J129").Value = 1002.36365
$ws.Range("L129").Value = 3007.09095
$ws.Range("N129").Value = -13007.09095
$ws.Range("H137").Value = 3572946
$ws.Range("I137").Value = 4001103.5
$ws.Range("J137").Value = 4966.6665
$ws.Range("K137").Value = 12003310.5
$ws.Range("L137").Value = 14899.9995
$ws.Range("M137").Value = -12000760.5
$ws.Range("N137").Value = -19999.9995
$ws.Range("H138").Value = 2527856.8
$ws.Range("I138").Value = 1040.8334
$ws.Range("J138").Value = 3475412.8
$ws.Range("K138").Value = 3122.5002
$ws.Range("L138").Value = 10426238.4
$ws.Range("M138").Value = 2017.4998
$ws.Range("N138").Value = -10436518.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13876.7
$ws.Range("I32").Value = 14277.559
$ws.Range("K32").Value = 14277.559
$ws.Range("M32").Value = -13990.559
$ws.Range("H61").Value = 33400828
$ws.Range("I61").Value = 41708944
$ws.Range("J61").Value = 168355.17
$ws.Range("K61").Value = 41708944
$ws.Range("L61").Value = 168355.17
$ws.Range("M61").Value = -41708732
$ws.Range("N61").Value = -168779.17
$ws.Range("H74").Value = 4272505.5
$ws.Range("I74").Value = 5977209.5
$ws.Range("J74").Value = 60883.766
$ws.Range("K74").Value = 5977209.5
$ws.Range("L74").Value = 60883.766
$ws.Range("M74").Value = -5976335.5
$ws.Range("N74").Value = -62631.766
$ws.Range("H77").Value = 4272505.5
$ws.Range("I77").Value = 5977209.5
$ws.Range("J77").Value = 60883.766
$ws.Range("K77").Value = 29886047.5
$ws.Range("L77").Value = 304418.83
$ws.Range("M77").Value = -29881679.5
$ws.Range("N77").Value = -313154.83
$ws.Range("H102").Value = 28573248
$ws.Range("I102").Value = 35715684
$ws.Range("K102").Value = 35715684
$ws.Range("M102").Value = -35714062
$ws.Range("H132").Value = 53792.75
$ws.Range("I132").Value = 39412.73
$ws.Range("J132").Value = 80498.5
$ws.Range("K132").Value = 118238.19
$ws.Range("L132").Value = 241495.5
$ws.Range("M132").Value = -115708.19
$ws.Range("N132").Value = -246555.5
$ws.Range("H136").Value = 33400828
$ws.Range("I136").Value = 41708944
$ws.Range("J136").Value = 168355.17
$ws.Range("K136").Value = 125126832
$ws.Range("L136").Value = 505065.51
$ws.Range("M136").Value = -125124282
$ws.Range("N136").Value = -510165.51
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 27500
$ws.Range("J9").Value = 27500
$ws.Range("L9").Value = 27500
$ws.Range("N9").Value = -27836
$ws.Range("H42").Value = 140000
$ws.Range("J42").Value = 140000
$ws.Range("L42").Value = 140000
$ws.Range("N42").Value = -140656
$ws.Range("H43").Value = 140000
$ws.Range("J43").Value = 140000
$ws.Range("L43").Value = 140000
$ws.Range("N43").Value = -140362
$ws.Range("H44").Value = 18000
$ws.Range("J44").Value = 18000
$ws.Range("L44").Value = 18000
$ws.Range("N44").Value = -18994
$ws.Range("H94").Value = 1012.5
$ws.Range("I94").Value = 1100
$ws.Range("J94").Value = 960
$ws.Range("K94").Value = 1100
$ws.Range("L94").Value = 960
$ws.Range("M94").Value = -649
$ws.Range("N94").Value = -1862
$ws.Range("H99").Value = 1083.8846
$ws.Range("I99").Value = 1077.8948
$ws.Range("J99").Value = 1100.1428
$ws.Range("K99").Value = 1077.8948
$ws.Range("L99").Value = 1100.1428
$ws.Range("M99").Value = 420.1052
$ws.Range("N99").Value = -4096.1428
$ws.Range("H105").Value = 35715780
$ws.Range("I105").Value = 62501200
$ws.Range("J105").Value = 1883.3334
$ws.Range("K105").Value = 62501200
$ws.Range("L105").Value = 1883.3334
$ws.Range("M105").Value = -62499453
$ws.Range("N105").Value = -5377.3334
$ws.Range("H134").Value = 3861.1794
$ws.Range("I134").Value = 3092.276
$ws.Range("J134").Value = 6091
$ws.Range("K134").Value = 9276.828
$ws.Range("L134").Value = 18273
$ws.Range("M134").Value = -6741.828
$ws.Range("N134").Value = -23343
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2714.4614
$ws.Range("I31").Value = 1391.1482
$ws.Range("J31").Value = 5691.9165
$ws.Range("K31").Value = 1391.1482
$ws.Range("L31").Value = 5691.9165
$ws.Range("M31").Value = -1096.1482
$ws.Range("N31").Value = -6281.9165
$ws.Range("H34").Value = 2714.4614
$ws.Range("I34").Value = 1391.1482
$ws.Range("J34").Value = 5691.9165
$ws.Range("K34").Value = 1391.1482
$ws.Range("L34").Value = 5691.9165
$ws.Range("M34").Value = -1189.1482
$ws.Range("N34").Value = -6095.9165
$ws.Range("H132").Value = 24680.047
$ws.Range("I132").Value = 1232.9714
$ws.Range("J132").Value = 127261
$ws.Range("K132").Value = 3698.9142
$ws.Range("L132").Value = 381783
$ws.Range("M132").Value = -1168.9142
$ws.Range("N132").Value = -386843
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7145458.5
$ws.Range("J4").Value = 7145458.5
$ws.Range("L4").Value = 21436375.5
$ws.Range("N4").Value = -21436599.5
$ws.Range("H75").Value = 2500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 7500
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -9496
$ws.Range("H78").Value = 2500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 22500
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -32484
$ws.Range("H87").Value = 24530.914
$ws.Range("I87").Value = 14501.77
$ws.Range("J87").Value = 30457.227
$ws.Range("K87").Value = 43505.31
$ws.Range("L87").Value = 91371.681
$ws.Range("M87").Value = -42257.31
$ws.Range("N87").Value = -93867.681
$ws.Range("H90").Value = 24530.914
$ws.Range("I90").Value = 14501.77
$ws.Range("J90").Value = 30457.227
$ws.Range("K90").Value = 130515.93
$ws.Range("L90").Value = 274115.043
$ws.Range("M90").Value = -124275.93
$ws.Range("N90").Value = -286595.043
$ws.Range("H131").Value = 1198.9344
$ws.Range("J131").Value = 1248.8392
$ws.Range("L131").Value = 3746.5176
$ws.Range("N131").Value = -13826.5176
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 73449.14
$ws.Range("I132").Value = 57009.11
$ws.Range("J132").Value = 103041.2
$ws.Range("K132").Value = 171027.33
$ws.Range("L132").Value = 309123.6
$ws.Range("M132").Value = -168497.33
$ws.Range("N132").Value = -314183.6
$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 105000
$ws.Range("N134").Value = -110070
$ws.Range("H141").Value = 16160
$ws.Range("I141").Value = 3390
$ws.Range("J141").Value = 23457.143
$ws.Range("K141").Value = 3390
$ws.Range("L141").Value = 23457.143
$ws.Range("M141").Value = 1790
$ws.Range("N141").Value = -33817.143
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 792
$ws.Range("I46").Value = 860
$ws.Range("J46").Value = 690
$ws.Range("K46").Value = 860
$ws.Range("L46").Value = 690
$ws.Range("M46").Value = -672
$ws.Range("N46").Value = -1066
$ws.Range("H132").Value = 36981.07
$ws.Range("I132").Value = 19674.393
$ws.Range("J132").Value = 69286.87
$ws.Range("K132").Value = 59023.179
$ws.Range("L132").Value = 207860.61
$ws.Range("M132").Value = -56493.179
$ws.Range("N132").Value = -212920.61
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 3250
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3250
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 3250
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -3804
$ws.Range("H56").Value = 31825
$ws.Range("J56").Value = 31825
$ws.Range("L56").Value = 31825
$ws.Range("N56").Value = -33253
$ws.Range("H136").Value = 28312.082
$ws.Range("I136").Value = 17920.17
$ws.Range("J136").Value = 69186.92999999999
$ws.Range("K136").Value = 53760.50999999999
$ws.Range("L136").Value = 207560.79
$ws.Range("M136").Value = -51210.50999999999
$ws.Range("N136").Value = -212660.79
